$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New data rows to append (Fecha, Pruebas Realizadas, Pruebas Positivas, Clinicamente Estables, Clinicamente Graves, Cuidados Intensivos)
$newRows = @(
    @(44189, 1244, 385, 487, 121, 23),
    @(44190, 968, 238, 468, 132, 25),
    @(44191, 1361, 438, 495, 126, 23)
)

$lastRow = 284
$startRow = $lastRow + 1

# Copy formatting from the last existing data row down into the new rows
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Range("A$lastRow`:F$lastRow").Copy()
    $ws.Range("A$r`:F$r").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}
$excel.CutCopyMode = $false

# Fill in the actual values
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
}

$newLastRow = $startRow + $newRows.Count - 1

# Resize the table to include the new rows
$tbl = $ws.ListObjects.Item("Condicion_Pacientes")
$newRange = $ws.Range("A1:F$newLastRow")
$tbl.Resize($newRange)

# Update sheet view: scroll position and selection
$ws.Range("D292").Select()
$excel.ActiveWindow.ScrollRow = 279
$excel.ActiveWindow.ScrollColumn = 1
